$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 116001
$ws.Range("B2").Value = 16001
$ws.Range("C2").Value = 16001
$ws.Range("D2").Value = "Lamina de hierro negro 3,17mm 1,22x2,44x1/8"
$ws.Range("E2").Value = 28561.95
$ws.Range("F2").Value = "unidad"
$ws.Range("G2").Value = "t"
$ws.Range("H2").Value = 13
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 16
$ws.Range("K2").Value = 7

# Row 3
$ws.Range("B3").Value = 16002
$ws.Range("C3").Value = 16002
$ws.Range("D3").Value = "TG Lamina Mold Tough 1/2 12mmx2,44 Gypsum Verde Resis Moho"
$ws.Range("E3").Value = 4986.7299999999996
$ws.Range("F3").Value = "unidad"
$ws.Range("G3").Value = "t"
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 16
$ws.Range("K3").Value = 7

# Row 4
$ws.Range("A4").Value = 6001
$ws.Range("B4").Value = 16003
$ws.Range("C4").Value = 16003
$ws.Range("D4").Value = "TG Lamina securock Glass 1/2x4x8 p/exterior aprox 30 kls"
$ws.Range("E4").Value = 11429.2
$ws.Range("F4").Value = "unidad"
$ws.Range("G4").Value = "t"
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 16
$ws.Range("K4").Value = 7

# Right-align F and G header + data (new style used)
$ws.Range("F1:G4").HorizontalAlignment = -4152
